$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last data row (row 10) into the new row 11 so the new
# record inherits the same number formats / borders / fonts as its
# neighbour before we overwrite the cell contents.
$ws.Range("A10:S10").Copy($ws.Range("A11:S11")) | Out-Null

# Fill in the new employee's data (order matches how the shared-string
# table is expected to grow: Name, Org code, Nationality, Team name,
# then the e-mail address last).
$ws.Range("A11").Value = "System Analyst"
$ws.Range("B11").Value = "Majeed, Jincy Karappamveettil"
$ws.Range("C11").Value = "LT52"
$ws.Range("J11").Value = "India"
$ws.Range("O11").Value = "LT52"
$ws.Range("P11").Value = "LT52 - Services & Call Center Team 123"
$ws.Range("D11").Value = "Jmajeed@kockw.com"
$ws.Range("E11").Value = 18450
$ws.Range("F11").Value = 11
$ws.Range("G11").Value = "F"
$ws.Range("H11").Value = 34
$ws.Range("I11").Value = ""
$ws.Range("K11").Value = "LT01"
$ws.Range("L11").Value = "LT01-CORPORATE INFORMATION TECHNOLOGY GROUP."
$ws.Range("M11").Value = "LX01"
$ws.Range("N11").Value = "LX01-(Planning & Innovation) Directorate"
$ws.Range("Q11").Value = 0
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 0

# The e-mail column is a mailto hyperlink, just like every row above it.
$ws.Hyperlinks.Add($ws.Range("D11"), "mailto:Jmajeed@kockw.com", [Type]::Missing, [Type]::Missing, "Jmajeed@kockw.com") | Out-Null

# Match the cursor position recorded in the saved file.
$ws.Range("A11").Select() | Out-Null
